$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from "Through 2022-11-28" to "Through 2022-11-29"
$ws.Name = "Through 2022-11-29"

# Update the column header label to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 11-29)"

# Update November's 2022-column value (row 12)
$ws.Range("I12").Value = 114

# Update the Total row's 2022-column value (row 14)
$ws.Range("I14").Value = 1512
